$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Fornece dados Pessoais " text to "Fornecer dados Pessoais "
$ws.Range("C9").Value = "Fornecer dados Pessoais "

# Update the exception text in A14 to reflect the new wording
$ws.Range("A14").Value = "Excepção 1               (passo 2)" + [char]10 + "[Cliente já se encontra registado no sistema]"

# Row 14 height changes from 60 to 90
$ws.Rows.Item(14).RowHeight = 90

# Update selection to B6
$ws.Range("B6").Select()
